$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds serial date values that need to be bumped
# from 45181 to 45182 for every data row (rows 2 through 115).
$ws.Range("C2:C115").Value = 45182
